$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-01-29 Monday" "2024-01-30 Tuesday"

Replace-Text "152÷3=" "424÷4="
Replace-Text "884÷7=" "369÷2="
Replace-Text "163÷6=" "458÷2="
Replace-Text "546÷9=" "711÷4="
Replace-Text "916÷5=" "307÷7="

Replace-Text "875÷9=" "455÷9="
Replace-Text "461÷7=" "176÷7="
Replace-Text "640÷8=" "898÷7="
Replace-Text "578÷5=" "565÷6="
Replace-Text "342÷4=" "746÷9="

Replace-Text "420÷9=" "360÷2="
Replace-Text "156÷6=" "918÷7="
Replace-Text "519÷4=" "456÷9="
Replace-Text "120÷3=" "382÷4="
Replace-Text "421÷6=" "900÷5="

Replace-Text "356÷4=" "205÷4="
Replace-Text "569÷5=" "345÷3="
Replace-Text "447÷2=" "438÷5="
Replace-Text "137÷8=" "155÷3="
Replace-Text "789÷7=" "470÷6="

Replace-Text "360÷8=" "501÷7="
Replace-Text "529÷8=" "372÷8="
Replace-Text "635÷6=" "627÷6="
Replace-Text "740÷7=" "457÷9="
Replace-Text "253÷7=" "206÷9="

Write-Host "Done"
